# Rename the embedded picture "name" metadata (wp:docPr / pic:cNvPr @name)
# for the three logo images in this document's headers/footers:
#
#   - Pearson logo in the "first page" footer  : image2.png -> image1.png
#   - Pearson logo in the "default" footer     : image2.png -> image1.png
#   - BTEC logo in the "first page" header     : image1.jpg -> image2.jpg
#
# InlineShape has no writable .Name property in the Word object model,
# so the standard technique is used: convert the inline picture to a
# (floating) Shape, rename it, then convert it back to an inline shape.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlinePicture($headerFooter, $newName) {
    if ($headerFooter.Exists -and $headerFooter.Range.InlineShapes.Count -gt 0) {
        $inlineShape = $headerFooter.Range.InlineShapes(1)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

# Footers: index 1 = default ("image2.png" / docPr id=2), index 2 = first page ("image2.png" / docPr id=3)
Rename-InlinePicture $sec.Footers(1) "image1.png"
Rename-InlinePicture $sec.Footers(2) "image1.png"

# Headers: index 2 = first page (BTec logo "image1.jpg" / docPr id=1)
Rename-InlinePicture $sec.Headers(2) "image2.jpg"

Write-Output "done"
